$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the visitor record (row 2) ---
# New date of visit
$ws.Range("A2").Value = "15/6/2025"
# New visitor name
$ws.Range("B2").Value = "Yuriana Montserrat Ibarra Granados"
# New control number (keep as text, like the original value)
$ws.Range("C2").Value = "'21420209"

# The visit type changes from "Revision Tesina" (G) to "Sala de Computacion" (I)
$ws.Range("G2").Value = ""
$ws.Range("I2").Value = "X"

# --- Page setup tweaks ---
# Keep fit-to-page settings explicit (fit to 1 page wide/tall)
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
